$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values to update (automatic electricity price update)
$ws.Range("A2").Value = 46071

$ws.Range("B2").Value = 15.61
$ws.Range("C2").Value = 15.31
$ws.Range("D2").Value = 5.31
$ws.Range("E2").Value = 3.89
$ws.Range("F2").Value = 3.72
$ws.Range("G2").Value = 3.53
$ws.Range("H2").Value = 2.8
$ws.Range("I2").Value = 3.55
$ws.Range("J2").Value = 7.28
$ws.Range("K2").Value = 16.13
$ws.Range("L2").Value = 5.66
$ws.Range("M2").Value = 1.67
$ws.Range("N2").Value = 0.42
$ws.Range("O2").Value = 0.28
$ws.Range("P2").Value = 0.24
$ws.Range("Q2").Value = 0.29
$ws.Range("R2").Value = 0.4
$ws.Range("S2").Value = 2.55
$ws.Range("T2").Value = 3.14
$ws.Range("U2").Value = 4.04
$ws.Range("V2").Value = 12.31
$ws.Range("W2").Value = 9.029999999999999
$ws.Range("X2").Value = 4.09
$ws.Range("Y2").Value = 1.39
$ws.Range("Z2").Value = 5.11

$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 10.03
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 15.46
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 11.7
$ws.Range("AG2").Value = "3h-23h"
